# "Actualización automatica" — refresh the column-header / dimension-code
# shared strings in the metadata sheet:
#   * a handful of accented header labels got mangled into their
#     UTF-8-bytes-read-as-cp1257 ("mojibake") form, and
#   * the "NNporcentaje" tokens in the slug/measure columns were
#     shortened to plain "NN".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 — accented header labels re-encoded as mojibake.
$ws.Range("D1").Value = "CĆ³digo de la provincia"
$ws.Range("G1").Value = "NĆŗmero total de explotaciones"
$ws.Range("K1").Value = "CĆ³digo del municipio"
$ws.Range("N1").Value = "CĆ³digo de la comarca"
$ws.Range("P1").Value = "AĆ±o"

# Row 2 — slug column, "NNporcentaje" -> "NN".
$ws.Range("A2").Value = "explotaciones-con-tierras-con-sau-75-y-100-de-su-propiedad"
$ws.Range("B2").Value = "explotaciones-con-tierras-con-sau-50-y-75-de-su-propiedad"
$ws.Range("L2").Value = "explotaciones-con-tierras-con-sau--25-de-su-propiedad"
$ws.Range("M2").Value = "explotaciones-con-tierras-con-sau-25-y-50-de-su-propiedad"

# Row 3 — "iaest-measure:" prefixed variant of the same slugs.
$ws.Range("A3").Value = "iaest-measure:explotaciones-con-tierras-con-sau-75-y-100-de-su-propiedad"
$ws.Range("B3").Value = "iaest-measure:explotaciones-con-tierras-con-sau-50-y-75-de-su-propiedad"
$ws.Range("L3").Value = "iaest-measure:explotaciones-con-tierras-con-sau--25-de-su-propiedad"
$ws.Range("M3").Value = "iaest-measure:explotaciones-con-tierras-con-sau-25-y-50-de-su-propiedad"
